# Data extraction for the acclimation set (bernhardt_2018, Figure 2 / Tetraselmis tetrahele
# population growth rate) and a screening-file scroll/selection refresh after the
# scopus-search pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New data rows 159-164 (bernhardt_2018, Figure 2, Tetraselmis tetrahele)
# ---------------------------------------------------------------------------

$rows = @(
    @{ Row=159; I=5;  J=0;  K=10; M=10; P=0.28682171000000001;  Q=0.33247422999999998;  R=0.18143124999999999;  S=0.12064242999999999 },
    @{ Row=160; I=10; J=5;  K=15; M=10; P=0.58914728999999999;  Q=0.71907215999999996;  R=0.048381670000000002; S=0.096513940000000006 },
    @{ Row=161; I=15; J=10; K=20; M=10; P=1.0465116299999999;   Q=0.88530927999999998;  R=0.13304958;           S=0.084449700000000003 },
    @{ Row=162; I=20; J=15; K=25; M=10; P=1.3139534900000001;   Q=1.2564432999999999;   R=0.25400373999999998;  S=0.31970242999999998 },
    @{ Row=163; I=24; J=19; K=29; M=10; P=1.4689922500000001;   Q=1.1520618600000001;   R=0.2237652;             S=0.28954182000000001 },
    @{ Row=164; I=27; J=22; K=32; M=10; P=1.5503876000000001;   Q=0.80412371000000005;  R=0.090715619999999997; S=0.69972606999999998 }
)

# Note: the shared-string table assigns new indices in first-seen order, so
# the cells below are touched in the same left-to-right/first-seen order the
# original authoring session produced (D, C, N, O, X, Y) to land each new
# string at the index the target workbook expects.

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = 1              # A study_id
    $ws.Cells.Item($row, 2).Value = 1              # B experiment_id

    $ws.Cells.Item($row, 4).Value = "Figure 2"         # D origin
    $ws.Cells.Item($row, 3).Value = "bernhardt_2018"   # C response_id

    $ws.Cells.Item($row, 5).Value = 1              # E org_level
    $ws.Cells.Item($row, 6).Value = 1              # F flux_treatment
    $ws.Cells.Item($row, 7).Value = 1              # G flux_pattern
    $ws.Cells.Item($row, 8).Value = 16             # H temp_ini

    # I mean_temp_constant - custom "0.0" number format + dark font
    $ws.Cells.Item($row, 9).Font.Color = 0
    $ws.Cells.Item($row, 9).NumberFormat = "0.0"
    $ws.Cells.Item($row, 9).Value = $r.I

    $ws.Cells.Item($row, 10).Value = $r.J           # J min_temp
    $ws.Cells.Item($row, 11).Value = $r.K           # K max_temp

    $ws.Cells.Item($row, 12).Formula = "=K" + $row + "-J" + $row   # L flux_range

    $ws.Cells.Item($row, 13).Value = $r.M           # M period_flux

    $ws.Cells.Item($row, 14).Value = "population growth rate "            # N resp_def
    $ws.Cells.Item($row, 15).Value = "population growth rate day^-1"      # O resp_units

    # P/Q constant_resp / flux_resp - dark font, default number format
    $ws.Cells.Item($row, 16).Font.Color = 0
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Font.Color = 0
    $ws.Cells.Item($row, 17).Value = $r.Q

    $ws.Cells.Item($row, 18).Value = $r.R           # R constant_variance
    $ws.Cells.Item($row, 19).Value = $r.S           # S flux_variance
    $ws.Cells.Item($row, 20).Value = 9              # T constant_samp
    $ws.Cells.Item($row, 21).Value = 9              # U flux_samp
    $ws.Cells.Item($row, 22).Value = 1              # V resp_quality
    $ws.Cells.Item($row, 23).Value = 1              # W variance_type

    $ws.Cells.Item($row, 24).Value = "Tetraselmis"  # X genus
    $ws.Cells.Item($row, 25).Value = "tetrahele"    # Y species

    $ws.Cells.Item($row, 26).Value = 1              # Z larger_group
    $ws.Cells.Item($row, 27).Value = 1              # AA exp_age
    $ws.Cells.Item($row, 28).Value = 1              # AB size

    $ws.Cells.Item($row, 36).Value = "y"            # AJ same_mean (y/n)
}

# ---------------------------------------------------------------------------
# Screening-file scroll/selection refresh after the scopus search pass
# ---------------------------------------------------------------------------

$ws.Range("AF169").Select()
